# Append new job postings and refresh timestamps (2025-11-21 18:23 JST)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: remove existing hyperlinks (reverse order avoids index-shift issues) ---
$existingLinks = @()
foreach ($hl in $ws.Hyperlinks) { $existingLinks += $hl }
for ($i = $existingLinks.Count - 1; $i -ge 0; $i--) { $existingLinks[$i].Delete() }

# --- Step 2: write the full A2:H14 data block (old rows pushed down, new rows inserted) ---
$ws.Range("A2").Value2 = "2025-11-21 18:23:52"
$ws.Range("B2").Value2 = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Range("C2").Value2 = "システム開発"
$ws.Range("D2").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value2 = "期限情報なし"
$ws.Range("F2").Value2 = "https://www.lancers.jp/work/detail/5434128"
$ws.Range("G2").Value2 = 368
$ws.Range("H2").Value2 = "🔥AI,Ai ◆開発"

$ws.Range("A3").Value2 = "2025-11-21 18:23:52"
$ws.Range("B3").Value2 = "企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)"
$ws.Range("C3").Value2 = "システム開発"
$ws.Range("D3").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value2 = "期限情報なし"
$ws.Range("F3").Value2 = "https://www.lancers.jp/work/detail/5434363"
$ws.Range("G3").Value2 = 348
$ws.Range("H3").Value2 = "🔥AI,Ai ◆コンサル"

$ws.Range("A4").Value2 = "2025-11-21 18:23:52"
$ws.Range("B4").Value2 = "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"
$ws.Range("C4").Value2 = "システム開発"
$ws.Range("D4").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value2 = "期限情報なし"
$ws.Range("F4").Value2 = "https://www.lancers.jp/work/detail/5439158"
$ws.Range("G4").Value2 = 303
$ws.Range("H4").Value2 = "🔥AI,Ai"

$ws.Range("A5").Value2 = "2025-11-21 18:23:52"
$ws.Range("B5").Value2 = "製造業のR&D支援!「プロセスデータ解析」「音響異常検知」のAIエンジニア募集"
$ws.Range("C5").Value2 = "システム開発"
$ws.Range("D5").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E5").Value2 = "期限情報なし"
$ws.Range("F5").Value2 = "https://www.lancers.jp/work/detail/5439165"
$ws.Range("G5").Value2 = 303
$ws.Range("H5").Value2 = "🔥AI,Ai"

$ws.Range("A6").Value2 = "2025-11-21 18:23:52"
$ws.Range("B6").Value2 = "【謝礼2,000円】AIに興味のあるエンジニアの方へ|45分だけお話を聞かせてください"
$ws.Range("C6").Value2 = "システム開発"
$ws.Range("D6").Value2 = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E6").Value2 = "期限情報なし"
$ws.Range("F6").Value2 = "https://www.lancers.jp/work/detail/5438171"
$ws.Range("G6").Value2 = 295
$ws.Range("H6").Value2 = "🔥AI,Ai"

$ws.Range("A7").Value2 = "2025-11-21 18:23:52"
$ws.Range("B7").Value2 = "初回 Pythonのテキストエディターに機能追加依頼"
$ws.Range("C7").Value2 = "システム開発"
$ws.Range("D7").Value2 = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E7").Value2 = "期限情報なし"
$ws.Range("F7").Value2 = "https://www.lancers.jp/work/detail/5439127"
$ws.Range("G7").Value2 = 190
$ws.Range("H7").Value2 = "🔥Python"

$ws.Range("A8").Value2 = "2025-11-21 18:23:52"
$ws.Range("B8").Value2 = "急募 【急募】MT4/MT5用FX自動売買システムの開発者募集"
$ws.Range("C8").Value2 = "システム開発"
$ws.Range("D8").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E8").Value2 = "期限情報なし"
$ws.Range("F8").Value2 = "https://www.lancers.jp/work/detail/5439201"
$ws.Range("G8").Value2 = 83
$ws.Range("H8").Value2 = "◆開発"

$ws.Range("A9").Value2 = "2025-11-21 18:23:52"
$ws.Range("B9").Value2 = "【データベース化】エクセル管理台帳の視覚化と検索機能強化"
$ws.Range("C9").Value2 = "システム開発"
$ws.Range("D9").Value2 = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E9").Value2 = "期限情報なし"
$ws.Range("F9").Value2 = "https://www.lancers.jp/work/detail/5438740"
$ws.Range("G9").Value2 = 30
$ws.Range("H9").Value2 = "◇管理"

$ws.Range("A10").Value2 = "2025-11-21 18:23:52"
$ws.Range("B10").Value2 = "限定公開 限定公開の仕事"
$ws.Range("C10").Value2 = "システム開発"
$ws.Range("D10").Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E10").Value2 = "期限情報なし"
$ws.Range("F10").Value2 = "https://www.lancers.jp/work/detail/5439193"
$ws.Range("G10").Value2 = 25
$ws.Range("H10").ClearContents()

$ws.Range("A11").Value2 = "2025-11-21 18:23:52"
$ws.Range("B11").Value2 = "【長期案件あり】Microsoft 365 / Azure / AWS エンジニアを急募!"
$ws.Range("C11").Value2 = "システム開発"
$ws.Range("D11").Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E11").Value2 = "期限情報なし"
$ws.Range("F11").Value2 = "https://www.lancers.jp/work/detail/5439182"
$ws.Range("G11").Value2 = 25
$ws.Range("H11").ClearContents()

$ws.Range("A12").Value2 = "2025-11-21 18:23:52"
$ws.Range("B12").Value2 = "【教育分野】新プロジェクトのPM募集!企画整理とチーム構築"
$ws.Range("C12").Value2 = "システム開発"
$ws.Range("D12").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E12").Value2 = "期限情報なし"
$ws.Range("F12").Value2 = "https://www.lancers.jp/work/detail/5438369"
$ws.Range("G12").Value2 = 18
$ws.Range("H12").ClearContents()

$ws.Range("A13").Value2 = "2025-11-21 18:23:52"
$ws.Range("B13").Value2 = "【急募】instagramとSTORES連携で商品販売を実現したい"
$ws.Range("C13").Value2 = "システム開発"
$ws.Range("D13").Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E13").Value2 = "期限情報なし"
$ws.Range("F13").Value2 = "https://www.lancers.jp/work/detail/5438567"
$ws.Range("G13").Value2 = 13
$ws.Range("H13").ClearContents()

$ws.Range("A14").Value2 = "2025-11-21 18:23:52"
$ws.Range("B14").Value2 = "PowerAutomate でWorepress記事を自動作成"
$ws.Range("C14").Value2 = "システム開発"
$ws.Range("D14").Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E14").Value2 = "期限情報なし"
$ws.Range("F14").Value2 = "https://www.lancers.jp/work/detail/5438092"
$ws.Range("G14").Value2 = 13
$ws.Range("H14").ClearContents()

# --- Step 3: (re)create hyperlinks on the URL column, then restore the shared Hyperlink style ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5434128") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5434363") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5439158") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5439165") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5438171") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5439127") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5439201") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5438740") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5439193") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5439182") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5438369") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5438567") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5438092") | Out-Null

$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("F14").Style = "Hyperlink"

# --- Step 4: column width changes (B: 45 -> 51, H: 12 -> 14 chars) ---
$ws.Columns.Item(2).ColumnWidth = 50.125
$ws.Columns.Item(8).ColumnWidth = 13.125

Write-Output "done"